$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

$ws.Cells.Item($row, 1).Value = 94973710
$ws.Cells.Item($row, 2).Value = 56543
$ws.Cells.Item($row, 3).Value = "Ovaliderad"
$ws.Cells.Item($row, 4).Value = "NT"
$ws.Cells.Item($row, 5).Value = 103021
$ws.Cells.Item($row, 6).Value = "Talltita"
$ws.Cells.Item($row, 7).Value = "Poecile montanus"
$ws.Cells.Item($row, 8).Value = "(Conrad von Baldenstein, 1827)"

$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 9).Value = "4"

$ws.Cells.Item($row, 11).Formula = '=""'
$ws.Cells.Item($row, 12).Formula = '=""'

$ws.Cells.Item($row, 13).Value = "pulli/nyligen flygga ungar"

$ws.Cells.Item($row, 14).Formula = '=""'

$ws.Cells.Item($row, 16).Value = "S Bådamossen A 32329-2021, Hl"
$ws.Cells.Item($row, 17).Value = 332576
$ws.Cells.Item($row, 18).Value = 6380968
$ws.Cells.Item($row, 19).Value = 100
$ws.Cells.Item($row, 20).Value = "Halland"
$ws.Cells.Item($row, 21).Value = "Kungsbacka"
$ws.Cells.Item($row, 22).Value = "Halland"
$ws.Cells.Item($row, 23).Value = "Tölö"

$ws.Cells.Item($row, 25).NumberFormat = "@"
$ws.Cells.Item($row, 25).Value = "2021-07-19"

$ws.Cells.Item($row, 27).NumberFormat = "@"
$ws.Cells.Item($row, 27).Value = "2021-07-19"

$ws.Cells.Item($row, 29).Value = "Två vuxna minst 2 juvenila"

$ws.Cells.Item($row, 30).Value = $false
$ws.Cells.Item($row, 31).Value = $false
$ws.Cells.Item($row, 33).Value = $false

$ws.Cells.Item($row, 46).Formula = '=""'

$ws.Cells.Item($row, 49).Value = "Christer Johansson"
$ws.Cells.Item($row, 50).Value = "Christer Johansson"

$ws.Cells.Item($row, 51).Formula = '=""'
